$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the log with two new rows, copying the formatting of the last
# existing data row (row 5) down to rows 6 and 7.
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)
$ws.Range("A7:C7").PasteSpecial(-4122)

$ws.Range("A6").Value = 43888
$ws.Range("B6").Value = 1.5
$ws.Range("C6").Value = "J'ai revu le MCD avec M. Benzonana et M. Konutse"

$ws.Range("A7").Value = 43889
$ws.Range("B7").Value = 1.5
$ws.Range("C7").Value = "J'ai refait le MCD et MLD avec l'aide de M. Benzonana et j'ai pu faire un retour sur l'avancement du projet"

# Wrap the (now longer) description text for the whole column.
$ws.Range("C2:C7").WrapText = $true
$ws.Rows.Item(7).AutoFit()

$ws.Range("C14").Select()
